# Updates cryptos list price/volume figures (and restores the Maker/Aave
# row ordering) to match the latest scrape.
#
# Note: several "Price" values look like plain decimals (e.g. "228.99").
# Assigning such a string straight to .Value lets Excel auto-convert it to
# a numeric cell, which would silently change its stored type/formatting.
# Prefixing those specific values with a leading apostrophe forces Excel to
# keep them as literal text (exactly like typing '228.99 into a cell),
# matching the original text-only data model of this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.771.95"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").Value = "2.093.14"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'228.99"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D6").Value = "'0.618"
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("D7").Value = "'61.25"
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("D10").Value = "'0.0847"
$ws.Range("E10").Value = "  +0.32%  "
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").Value = "'15.34"
$ws.Range("E12").Value = "  +4.57%  "
$ws.Range("D13").Value = "2.402.84"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").Value = "'22.09"
$ws.Range("E14").Value = "  -1.82%  "
$ws.Range("D15").Value = "'0.805"
$ws.Range("E15").Value = "  +4.05%  "
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "2.088.35"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").Value = "38.698.34"
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("D19").Value = "'71.82"
$ws.Range("E19").Value = "  +2.09%  "
$ws.Range("E20").Value = "  +1.66%  "
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("D22").Value = "'228.03"
$ws.Range("E22").Value = "  +1.59%  "
$ws.Range("E24").Value = "  -2.63%  "
$ws.Range("D25").Value = "'2.35"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("D26").Value = "'171.58"
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("D27").Value = "'9.54"
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("D28").Value = "'0.138"
$ws.Range("E28").Value = "  +4.78%  "
$ws.Range("D29").Value = "'1.42"
$ws.Range("E29").Value = "  +5.49%  "
$ws.Range("D30").Value = "'19.34"
$ws.Range("E30").Value = "  +1.65%  "
$ws.Range("E31").Value = "  +2.95%  "
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("E33").Value = "  +2.04%  "
$ws.Range("D34").Value = "'4.74"
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("D35").Value = "'0.0621"
$ws.Range("E35").Value = "  +2.51%  "
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("D37").Value = "'2.40"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("D39").Value = "'0.998"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").Value = "'18.22"
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("E41").Value = "  +4.24%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.536.84"
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'101.01"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("E44").Value = "  -1.20%  "
$ws.Range("D45").Value = "'0.0910"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("D47").Value = "'7.67"
$ws.Range("E47").Value = "  +5.59%  "
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("E50").Value = "  -1.36%  "
$ws.Range("D51").Value = "2.290.23"
$ws.Range("E51").Value = "  -0.14%  "
